# "Additional companies sent for questionaire"
# The "Parent company" and "Location County/City" columns are removed from
# the known-locomotive list, shifting the remaining columns left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column E ("Location County/City") first so column B's position is
# unaffected by the later deletion.
$ws.Range("E:E").Delete()
# Delete column B ("Parent company").
$ws.Range("B:B").Delete()

# Match the author's final active-cell selection.
[void]$ws.Range("G6").Select()
